$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 2.62
$ws.Range("J2").Value = 3.1
$ws.Range("K2").Value = 2.1
$ws.Range("L2").Value = 3.15
$ws.Range("M2").Value = 9.5
$ws.Range("P2").Value = 3.05
$ws.Range("Q2").Value = 1.88
$ws.Range("U2").Value = 1.65
$ws.Range("V2").Value = 2
$ws.Range("AD2").Value = 6.1
$ws.Range("AH2").Value = 8.75
$ws.Range("AJ2").Value = 9.5
$ws.Range("AK2").Value = 30
$ws.Range("AM2").Value = 28
$ws.Range("AN2").Value = 4.55
$ws.Range("AO2").Value = 13.5
$ws.Range("AP2").Value = 19
$ws.Range("AQ2").Value = 55
$ws.Range("AR2").Value = 80
$ws.Range("AS2").Value = 200
$ws.Range("AT2").Value = 2.65
$ws.Range("AU2").Value = 6.4
$ws.Range("AV2").Value = 50
$ws.Range("AW2").Value = 4.65
$ws.Range("AX2").Value = 13.5
$ws.Range("AY2").Value = 19
$ws.Range("AZ2").Value = 60
$ws.Range("BA2").Value = 80
$ws.Range("BB2").Value = 200
$ws.Range("M3").Value = 1.06
$ws.Range("N3").Value = 5.8
$ws.Range("G5").Value = 4.75
$ws.Range("H5").Value = 3.5
$ws.Range("I5").Value = 1.62
$ws.Range("L5").Value = 2.3
$ws.Range("M5").Value = 1.06
$ws.Range("N5").Value = 9.5
$ws.Range("AC5").Value = 9.5
$ws.Range("AD5").Value = 7
$ws.Range("AL5").Value = 15
$ws.Range("AN5").Value = 6.5
$ws.Range("AW5").Value = 3.6
$ws.Range("AZ5").Value = 29
$ws.Range("Q6").Value = 1.85
$ws.Range("R6").Value = 1.95
$ws.Range("G7").Value = 3.1
$ws.Range("I7").Value = 2.25
$ws.Range("L7").Value = 3.1
$ws.Range("N7").Value = 7.5
$ws.Range("Y7").Value = 12
$ws.Range("AW7").Value = 4.33
$ws.Range("AZ7").Value = 51
$ws.Range("BA7").Value = 81
$ws.Range("K9").Value = 2.1
$ws.Range("M9").Value = 1.06
$ws.Range("N9").Value = 10
$ws.Range("X9").Value = 9.5
$ws.Range("AB9").Value = 29
$ws.Range("AC9").Value = 10
$ws.Range("AY9").Value = 29

Write-Output "Applied odds updates to rows 2, 3, 5, 6, 7, 9"
